$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 259.92307
$ws.Range("I28").Value = 256.45456
$ws.Range("K28").Value = 256.45456
$ws.Range("M28").Value = 228.54544
$ws.Range("H41").Value = 584.6667
$ws.Range("I41").Value = 288.125
$ws.Range("J41").Value = 923.5714
$ws.Range("K41").Value = 288.125
$ws.Range("L41").Value = 923.5714
$ws.Range("M41").Value = 151.875
$ws.Range("N41").Value = -1803.5714
$ws.Range("H58").Value = 2156
$ws.Range("J58").Value = 3483.3333
$ws.Range("L58").Value = 10449.9999
$ws.Range("N58").Value = -10749.9999
$ws.Range("H75").Value = 60000
$ws.Range("J75").Value = 60000
$ws.Range("L75").Value = 60000
$ws.Range("N75").Value = -61872
$ws.Range("H78").Value = 60000
$ws.Range("J78").Value = 60000
$ws.Range("L78").Value = 180000
$ws.Range("N78").Value = -189360
$ws.Range("H107").Value = 1239.9062
$ws.Range("I107").Value = 1228.6086
$ws.Range("J107").Value = 1268.7778
$ws.Range("K107").Value = 1228.6086
$ws.Range("L107").Value = 1268.7778
$ws.Range("M107").Value = 691.3914
$ws.Range("N107").Value = -5108.7778
$ws.Range("H132").Value = 2647.6858
$ws.Range("I132").Value = 2627.1562
$ws.Range("K132").Value = 7881.4686
$ws.Range("M132").Value = -5351.4686
$ws.Range("H138").Value = 2536.9841
$ws.Range("I138").Value = 890.6667
$ws.Range("J138").Value = 3195.5112
$ws.Range("K138").Value = 2672.0001
$ws.Range("L138").Value = 9586.533599999999
$ws.Range("M138").Value = 2467.9999
$ws.Range("N138").Value = -19866.5336

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8095.1772
$ws.Range("I32").Value = 4988.727
$ws.Range("K32").Value = 4988.727
$ws.Range("M32").Value = -4701.727
$ws.Range("H45").Value = 2669
$ws.Range("J45").Value = 2714
$ws.Range("L45").Value = 2714
$ws.Range("N45").Value = -3468
$ws.Range("H61").Value = 3265.8333
$ws.Range("I61").Value = 3098.8462
$ws.Range("K61").Value = 3098.8462
$ws.Range("M61").Value = -2886.8462
$ws.Range("H88").Value = 2020.7858
$ws.Range("I88").Value = 2009.0588
$ws.Range("J88").Value = 2038.909
$ws.Range("K88").Value = 2009.0588
$ws.Range("L88").Value = 2038.909
$ws.Range("M88").Value = -1603.0588
$ws.Range("N88").Value = -2850.909
$ws.Range("H91").Value = 2020.7858
$ws.Range("I91").Value = 2009.0588
$ws.Range("J91").Value = 2038.909
$ws.Range("K91").Value = 2009.0588
$ws.Range("L91").Value = 2038.909
$ws.Range("M91").Value = -605.0588
$ws.Range("N91").Value = -4846.909
$ws.Range("H122").Value = 1788.0454
$ws.Range("I122").Value = 1701.9474
$ws.Range("K122").Value = 5105.8422
$ws.Range("M122").Value = -2655.8422
$ws.Range("H132").Value = 2228.3257
$ws.Range("I132").Value = 1921.05
$ws.Range("J132").Value = 2495.5217
$ws.Range("K132").Value = 5763.15
$ws.Range("L132").Value = 7486.5651
$ws.Range("M132").Value = -3233.15
$ws.Range("N132").Value = -12546.5651
$ws.Range("H136").Value = 3265.8333
$ws.Range("I136").Value = 3098.8462
$ws.Range("K136").Value = 9296.5386
$ws.Range("M136").Value = -6746.5386

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 42797
$ws.Range("I81").Value = 50000
$ws.Range("K81").Value = 50000
$ws.Range("M81").Value = -48939
$ws.Range("H84").Value = 42797
$ws.Range("I84").Value = 50000
$ws.Range("K84").Value = 150000
$ws.Range("M84").Value = -144696
$ws.Range("H134").Value = 6025.05
$ws.Range("I134").Value = 6025.05
$ws.Range("K134").Value = 18075.15
$ws.Range("M134").Value = -15540.15

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2536.125
$ws.Range("I16").Value = 1445.75
$ws.Range("J16").Value = 3626.5
$ws.Range("K16").Value = 1445.75
$ws.Range("L16").Value = 3626.5
$ws.Range("M16").Value = -1158.75
$ws.Range("N16").Value = -4200.5
$ws.Range("H33").Value = 4281
$ws.Range("I33").Value = 4281
$ws.Range("K33").Value = 4281
$ws.Range("M33").Value = -3902
$ws.Range("H36").Value = 22249.666
$ws.Range("I36").Value = 22249.666
$ws.Range("K36").Value = 22249.666
$ws.Range("M36").Value = -21861.666
$ws.Range("H40").Value = 22249.666
$ws.Range("I40").Value = 22249.666
$ws.Range("K40").Value = 22249.666
$ws.Range("M40").Value = -22089.666
$ws.Range("H42").Value = 6332.6665
$ws.Range("I42").Value = 6332.6665
$ws.Range("K42").Value = 6332.6665
$ws.Range("M42").Value = -5739.6665
$ws.Range("H53").Value = 50000
$ws.Range("J53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("N53").Value = -51214
$ws.Range("H113").Value = 2536.125
$ws.Range("I113").Value = 1445.75
$ws.Range("J113").Value = 3626.5
$ws.Range("K113").Value = 1445.75
$ws.Range("L113").Value = 3626.5
$ws.Range("M113").Value = 724.25
$ws.Range("N113").Value = -7966.5
$ws.Range("H122").Value = 2107.4285
$ws.Range("I122").Value = 2272.6667
$ws.Range("J122").Value = 1611.7142
$ws.Range("K122").Value = 6818.000100000001
$ws.Range("L122").Value = 4835.142599999999
$ws.Range("M122").Value = -4368.000100000001
$ws.Range("N122").Value = -9735.142599999999
$ws.Range("H132").Value = 7392.7417
$ws.Range("I132").Value = 6513.2607
$ws.Range("J132").Value = 9921.25
$ws.Range("K132").Value = 19539.7821
$ws.Range("L132").Value = 29763.75
$ws.Range("M132").Value = -17009.7821
$ws.Range("N132").Value = -34823.75
$ws.Range("H134").Value = 3452.5386
$ws.Range("I134").Value = 3228.875
$ws.Range("J134").Value = 6136.5
$ws.Range("K134").Value = 9686.625
$ws.Range("L134").Value = 18409.5
$ws.Range("M134").Value = -7151.625
$ws.Range("N134").Value = -23479.5

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 4112.8887
$ws.Range("I12").Value = 508.6
$ws.Range("K12").Value = 1525.8
$ws.Range("M12").Value = -1352.8
$ws.Range("H141").Value = 3262.6365
$ws.Range("I141").Value = 1984.1428
$ws.Range("K141").Value = 5952.428400000001
$ws.Range("M141").Value = -772.4284000000007

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2905.9565
$ws.Range("I102").Value = 2992.2856
$ws.Range("K102").Value = 2992.2856
$ws.Range("M102").Value = -1370.2856
$ws.Range("H122").Value = 5530.913
$ws.Range("I122").Value = 4470.3125
$ws.Range("J122").Value = 7955.143
$ws.Range("K122").Value = 13410.9375
$ws.Range("L122").Value = 23865.429
$ws.Range("M122").Value = -10960.9375
$ws.Range("N122").Value = -28765.429
$ws.Range("H132").Value = 71801.39999999999
$ws.Range("I132").Value = 76501.64
$ws.Range("K132").Value = 229504.92
$ws.Range("M132").Value = -226974.92

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1523.6923
$ws.Range("J16").Value = 1485.2858
$ws.Range("L16").Value = 1485.2858
$ws.Range("N16").Value = -1825.2858
$ws.Range("H22").Value = 700
$ws.Range("I22").Value = 700
$ws.Range("K22").Value = 700
$ws.Range("M22").Value = -405
$ws.Range("H27").Value = 700
$ws.Range("I27").Value = 700
$ws.Range("K27").Value = 700
$ws.Range("M27").Value = -593
$ws.Range("H40").Value = 5582.6
$ws.Range("I40").Value = 5582.6
$ws.Range("K40").Value = 5582.6
$ws.Range("M40").Value = -5446.6
$ws.Range("H136").Value = 3018.4119
$ws.Range("I136").Value = 2083.0908
$ws.Range("K136").Value = 6249.2724
$ws.Range("M136").Value = -3699.2724

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3727.2856
$ws.Range("I122").Value = 3852.9092
$ws.Range("J122").Value = 3266.6667
$ws.Range("K122").Value = 11558.7276
$ws.Range("L122").Value = 9800.000100000001
$ws.Range("M122").Value = -9108.7276
$ws.Range("N122").Value = -14700.0001
$ws.Range("H126").Value = 3499.4285
$ws.Range("I126").Value = 3198.8
$ws.Range("J126").Value = 4251
$ws.Range("K126").Value = 9596.400000000001
$ws.Range("L126").Value = 12753
$ws.Range("M126").Value = -7126.400000000001
$ws.Range("N126").Value = -17693
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("H136").Value = 334981.3
$ws.Range("I136").Value = 386286.2
$ws.Range("J136").Value = 1499.5
$ws.Range("K136").Value = 1158858.6
$ws.Range("L136").Value = 4498.5
$ws.Range("M136").Value = -1156308.6
$ws.Range("N136").Value = -9598.5

# --- Clear cell on WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N132").ClearContents()
